$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-7 ("Anushka" block): Emp ID changed to a new Mongo-style id.
$ws.Range("A2").Value = "69067bfe21810f4deae03b5c"
$ws.Range("A3").Value = "69067bfe21810f4deae03b5c"
$ws.Range("A4").Value = "69067bfe21810f4deae03b5c"
$ws.Range("A5").Value = "69067bfe21810f4deae03b5c"
$ws.Range("A6").Value = "69067bfe21810f4deae03b5c"
$ws.Range("A7").Value = "69067bfe21810f4deae03b5c"

# Rows 14-19 ("Trisha" block): Emp ID changed to a new Mongo-style id.
$ws.Range("A14").Value = "6907a4d6f37984871bfd0e75"
$ws.Range("A15").Value = "6907a4d6f37984871bfd0e75"
$ws.Range("A16").Value = "6907a4d6f37984871bfd0e75"
$ws.Range("A17").Value = "6907a4d6f37984871bfd0e75"
$ws.Range("A18").Value = "6907a4d6f37984871bfd0e75"
$ws.Range("A19").Value = "6907a4d6f37984871bfd0e75"

# Trailing blank row: clear the stray formatted-but-empty B38 cell entirely.
$ws.Range("B38").Clear()

# Leave the cursor/selection where the author left off before saving.
[void]$ws.Range("A19").Select()
